$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 76 (the current "FRIDA" row) and insert the copy above it,
# preserving its exact formatting (borders, fills, merges), then shift the
# original row down to 77.
$ws.Rows.Item(76).Copy()
$ws.Rows.Item(76).Insert()

# Fill the newly inserted row 76 with the new item's data.
$ws.Cells.Item(76, 1).Value = 70
$ws.Cells.Item(76, 3).Value = "مجموعه لونا 3 قطع"
$ws.Cells.Item(76, 8).Value = "6:0"
$ws.Cells.Item(76, 12).Value = "0"
$ws.Cells.Item(76, 14).Value = "150.00"
$ws.Cells.Item(76, 16).Value = "150.0000"
$ws.Cells.Item(76, 17).Value = "1:0"

# Update the running subtotal to include the new item's price.
$ws.Cells.Item(78, 16).Value = 4572.4949999999999

# Refresh the "uploaded" timestamp footer.
$ws.Cells.Item(79, 1).Value = "Friday, 19 September, 2025 9:07 PM"
